$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a trailing space to the project names for "קיוסק ת"א" and
# "ריז 2 ת"א" rows (the user edited these values in column E, which is why
# Excel appends new shared-string entries and re-points the affected cells).
$ws.Range("E2").Value = "קיוסק ת""א "
$ws.Range("E3").Value = "קיוסק ת""א "
$ws.Range("E5").Value = "קיוסק ת""א "

$ws.Range("E6").Value = "ריז 2 ת""א "
$ws.Range("E7").Value = "ריז 2 ת""א "
$ws.Range("E9").Value = "ריז 2 ת""א "

# Scroll the sheet view back to the top and move the active selection,
# matching the saved view state (selection on K12, no frozen topLeftCell).
$ws.Range("K12").Select()
